$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -15.44523978092002
$ws.Range("F2").Value = -25.11195135116577
$ws.Range("I2").Value = 6.797258019447327
$ws.Range("J2").Value = -78.53172183036804

$ws.Range("B3").Value = -17.1990480123427
$ws.Range("F3").Value = -24.08938252925873
$ws.Range("I3").Value = 6.872270703315735
$ws.Range("J3").Value = -76.51601421833038

$ws.Range("B4").Value = -18.80406695003762
$ws.Range("F4").Value = -23.05302369594574
$ws.Range("I4").Value = 6.78860878944397
$ws.Range("J4").Value = -75.03032040596008

$ws.Range("B5").Value = -19.54645468242234
$ws.Range("F5").Value = -22.15524387359619
$ws.Range("I5").Value = 6.56968080997467
$ws.Range("J5").Value = -74.09579634666443

$ws.Range("B6").Value = -19.80004769015454
$ws.Range("F6").Value = -21.62751197814941
$ws.Range("I6").Value = 6.407713055610657
$ws.Range("J6").Value = -73.61052095890045

$ws.Range("B7").Value = -19.80004769015454
$ws.Range("F7").Value = -21.62751197814941
$ws.Range("I7").Value = 6.407713055610657
$ws.Range("J7").Value = -73.61052095890045

$ws.Range("B8").Value = -19.64176719690795
$ws.Range("F8").Value = -21.60710525512695
$ws.Range("I8").Value = 6.377492427825928
$ws.Range("J8").Value = -73.63119733333588

$ws.Range("B9").Value = -18.01807480916636
$ws.Range("F9").Value = -22.17591226100922
$ws.Range("I9").Value = 6.330997347831726
$ws.Range("J9").Value = -74.57762789726257

$ws.Range("B10").Value = -15.25805431506797
$ws.Range("F10").Value = -24.3661572933197
$ws.Range("I10").Value = 6.483354091644287
$ws.Range("J10").Value = -78.02318370342255

$ws.Range("B11").Value = -13.78136061189934
$ws.Range("F11").Value = -26.66655802726746
$ws.Range("I11").Value = 6.589540958404541
$ws.Range("J11").Value = -82.2586430311203

$ws.Range("B12").Value = -13.5017184536905
$ws.Range("F12").Value = -27.65814089775085
$ws.Range("I12").Value = 6.54293966293335
$ws.Range("J12").Value = -84.33668959140778

$ws.Range("B13").Value = -13.71766596963676
$ws.Range("F13").Value = -28.38467574119568
$ws.Range("I13").Value = 6.490947246551514
$ws.Range("J13").Value = -84.83413124084473

$ws.Range("B14").Value = -13.98665641243042
$ws.Range("F14").Value = -28.34670257568359
$ws.Range("I14").Value = 6.41607391834259
$ws.Range("J14").Value = -83.91497611999512

$ws.Range("B15").Value = -14.12430551401349
$ws.Range("F15").Value = -27.88725340366364
$ws.Range("I15").Value = 6.304742097854614
$ws.Range("J15").Value = -82.89813685417175

$ws.Range("B16").Value = -14.24435629204004
$ws.Range("F16").Value = -28.09894490242004
$ws.Range("I16").Value = 6.413556218147278
$ws.Range("J16").Value = -82.81659340858459

$ws.Range("B17").Value = -14.36836242720551
$ws.Range("F17").Value = -28.38948857784271
$ws.Range("I17").Value = 6.554604530334473
$ws.Range("J17").Value = -82.82983100414276

$ws.Range("B18").Value = -14.5920262637826
$ws.Range("F18").Value = -27.66289269924164
$ws.Range("I18").Value = 6.540995359420776
$ws.Range("J18").Value = -81.48527276515961

$ws.Range("B19").Value = -14.67099467912431
$ws.Range("F19").Value = -27.37504005432129
$ws.Range("I19").Value = 6.604041576385498
$ws.Range("J19").Value = -81.05589663982391

$ws.Range("B20").Value = -14.03366148806253
$ws.Range("F20").Value = -27.38587117195129
$ws.Range("I20").Value = 6.591673612594604
$ws.Range("J20").Value = -82.43313658237457

$ws.Range("B21").Value = -13.47106495284288
$ws.Range("F21").Value = -28.7396000623703
$ws.Range("I21").Value = 6.662915945053101
$ws.Range("J21").Value = -86.33434653282166

$ws.Range("B22").Value = -13.6531680320013
$ws.Range("F22").Value = -30.26678884029388
$ws.Range("I22").Value = 6.823388814926147
$ws.Range("J22").Value = -89.25043654441833

$ws.Range("B23").Value = -13.97864837094039
$ws.Range("F23").Value = -31.57401275634766
$ws.Range("I23").Value = 6.996102929115295
$ws.Range("J23").Value = -91.5533035993576

$ws.Range("B24").Value = -14.13984127349045
$ws.Range("F24").Value = -32.1827267408371
$ws.Range("I24").Value = 7.11737596988678
$ws.Range("J24").Value = -92.46497881412506

$ws.Range("B25").Value = -13.68327605135119
$ws.Range("F25").Value = -30.50069725513458
$ws.Range("I25").Value = 7.20659077167511
$ws.Range("J25").Value = -88.27713894844055

